$d = $word.ActiveDocument

# --- Change 1: "(4) Grab a chip (index finger's tip inside the chip)" ---
# becomes "(4) Grab a chip (index finger's tip inside the chip object)"
$d.Content.Find.Execute(
    "tip inside the chip)", $true, $false, $false, $false, $false,
    $true, 1, $false, "tip inside the chip object)", 2)

# --- Change 2: paragraph (6) text is replaced ---
# "(6) <TBD all other gestures and utilities>"
# becomes "(6) When showing the cards in left hand, snap (only thumb, index,
#          middle finger extended and thumb is close to middle finger) to use black magic"
$d.Content.Find.Execute(
    "(6) <TBD all other gestures and utilities>", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(6) When showing the cards in left hand, snap (only thumb, index, middle finger extended and thumb is close to middle finger) to use black magic",
    2)

# --- Change 3: insert new paragraphs (7), (8), (9) and two blank ones ---
# The trailing bookmark ("_GoBack") originally sat alone in the final blank
# paragraph; move it into the middle of the new paragraph (9), matching
# where Word last left the edit point.
$d.Bookmarks.Item("_GoBack").Delete()

# Locate paragraph (6) (the one we just edited above).
$p6 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*black magic*") {
        $p6 = $p
        break
    }
}

$p6.Range.InsertParagraphAfter()
$p7 = $p6.Next()
$p7.Range.InsertAfter("`t(7) Thumb down to taunt (text UI)")

$p7.Range.InsertParagraphAfter()
$p8 = $p7.Next()
$p8.Range.InsertAfter("`t(8) Clap to say gg ")

$p8.Range.InsertParagraphAfter()
$p9 = $p8.Next()
$p9.Range.InsertAfter("`t(9) Pistol gesture (finger and thumb extended) for threatening, rotate pistol along x axis to fire")

$p9.Range.InsertParagraphAfter()
$pBlank1 = $p9.Next()

$pBlank1.Range.InsertParagraphAfter()

# Re-seat the _GoBack bookmark inside paragraph (9), right after
# "(finger and t" (i.e. before "humb extended)...").
$full = $d.Content.Text
$marker = "(9) Pistol gesture (finger and t"
$idx = $full.IndexOf($marker)
$pos = $idx + $marker.Length
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
